$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (MRN 7506405 / 42628) -- remaining rows shift up
$ws.Rows("2:2").Delete()

# Delete the old last row (was row 8, now row 7 after the shift above)
$ws.Rows("7:7").Delete()

# Set selection to mirror the post-delete state recorded in the file
$ws.Range("A2:XFD2").Select()
